# Commit: "data process error message, region"
#
# The underlying data on the sheet is unchanged; this commit only touches
# workbook/sheet presentation metadata:
#   - the worksheet tab is renamed from the generic "Sheet1" to "br_vs"
#   - the author's last selection on the sheet moves to H24 (was G10)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: Sheet1 -> br_vs
$ws.Name = "br_vs"

# Move/save the active cell selection to H24 (matches <selection activeCell="H24" sqref="H24"/>)
$ws.Range("H24").Select()
